$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "bien supongo"
$ws.Range("A2").Value = "ajjasjassj"
$ws.Range("A3").Value = "Es muy"
$ws.Range("A4").Value = "lindo"
$ws.Range("A5").Value = "todo bien"
$ws.Range("A6").Value = "vos"
